# Update planner to local settings.
# Moves/edits the meal-plan entries and recalculated calorie totals so
# the workbook matches the "local settings" (days/columns) for this planner.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 ("New Meal"): clear the Tuesday (E3) entry - it has moved to row 6.
$ws.Range("E3").Value = ""

# Row 5 ("test1"): reorder the Thursday (G5) meal lines - Halibut now listed first.
$ws.Range("G5").Value = "HALIBUT (2.0 0Z)" + [char]10 + "CRAB (6.0 0Z)"

# Row 6 ("k k"): Monday (D6) gains an Ezekiel Bread line ahead of the Arby's entry,
# Wednesday (F6) now carries the Greek Yogurt/Apple combo (moved from E3),
# and Friday (H6) gets the same combo, lines reordered (Apple first).
$ws.Range("D6").Value = "EZEKIEL BREAD (1.0 SLICE)" + [char]10 + "ARBYS, roast beef sandwich, classic (3.0 sandwich)"
$ws.Range("F6").Value = "GREEK YOGURT (1.0 CUP)" + [char]10 + "APPLE (1.0 MEDIUM)"
$ws.Range("H6").Value = "APPLE (1.0 MEDIUM)" + [char]10 + "GREEK YOGURT (1.0 CUP)"

# Row 9 (Calories): update the day total text and move the 223-calorie value
# from Tuesday (E9) to Wednesday (F9) and Friday (H9), matching the meals above.
$ws.Range("D9").Value = "1,163.0"
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 223
$ws.Range("H9").Value = 223
